$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 129, shifting existing rows 129-152 down to 130-153.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new record.
# Columns A,B,C,E,F,G,H,I,N,Q,R carry the same values the row above (old row 129,
# now at 130) already has, while D,J,K,L,M,O,P get the new figures from the diff.
$ws.Cells.Item(129, 1).Value = 10
$ws.Cells.Item(129, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(129, 3).Value = "La Araucanía"
$ws.Cells.Item(129, 4).Value = 44798
$ws.Cells.Item(129, 5).Value = 9
$ws.Cells.Item(129, 6).Value = 100112031
$ws.Cells.Item(129, 7).Value = "Poroto verde"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 50
$ws.Cells.Item(129, 11).Value = 35000
$ws.Cells.Item(129, 12).Value = 35000
$ws.Cells.Item(129, 13).Value = 35000
$ws.Cells.Item(129, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(129, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(129, 16).Value = 1400
$ws.Cells.Item(129, 17).Value = 25
$ws.Cells.Item(129, 18).Value = "Hortaliza"
